$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Data")

# Row 3-4: replace hardcoded German labels with translation placeholders
$ws.Range("A3").Value = "{generiertAmTitle}"
$ws.Range("A4").Value = "{faelligAmTitle}"

# Row 6: table header labels -> translation placeholders
$ws.Range("A6").Value = "{institutionTitle}"
$ws.Range("B6").Value = "{nachnameTitle}"
$ws.Range("C6").Value = "{vornameTitle}"
$ws.Range("D6").Value = "{geburtsdatumTitle}"
$ws.Range("E6").Value = "{verfuegungTitle}"
$ws.Range("F6").Value = "{vonTitle}"
$ws.Range("G6").Value = "{bisTitle}"
$ws.Range("H6").Value = "{bgPensumTitle}"
$ws.Range("I6").Value = "{betragCHFTitle}"
$ws.Range("J6").Value = "{korrekturTitle}"
$ws.Range("K6").Value = "{zahlungIgnorierenTitle}"

# Update selection to match the saved workbook state (A6 selected)
$ws.Range("A6").Select() | Out-Null
